$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Пример")

# The only real content change: the "Сайт организации" value cell (B10)
# changed from "www.stat.kg" to "www.stat.gov.kg".
$ws.Range("B10").Value = "www.stat.gov.kg"

# View-state tweak captured in the diff: selection moved from B16 to B10.
$ws.Range("B10").Select()
